# "Inclusão de CASSIA SILVA" - append a new employee row (CPF / PIS / COLABORADOR)
# to the bottom of the "Planilha1" table, mirroring the format of the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Copy the formatting (borders/style) of the last existing data row (56) down
# into the new row (57) so the appended record matches the rest of the table.
$ws.Range("A56:C56").Copy()
$ws.Range("A57:C57").PasteSpecial(-4122)

# New record: CPF, PIS, COLABORADOR
$ws.Range("A57").Value = 52561493857
$ws.Range("B57").Value = 20440736743
$ws.Range("C57").Value = "CASSIA PEREIRA MARQUES DA SILVA"

# Reflect the manual scroll/zoom/selection state left behind after the edit.
$excel.ActiveWindow.Zoom = 160
[void]$ws.Range("A58").Select()
